# Update cryptocurrency price/volume snapshot values in the worksheet.
# (Auto-generated update mirroring a GitHub Actions data refresh.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '93.780.52'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.14%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.057.49'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.42%  '

# Row 4
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.62%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '603.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.21%  '

# Row 7
$ws.Range("E7").Value = '  -2.69%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.372'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -9.26%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.799'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.75%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.050.63'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.54%  '

# Row 12
$ws.Range("E12").Value = '  -4.18%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '93.354.59'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.03%  '

# Row 14
$ws.Range("E14").Value = '  -7.80%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '33.13'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.22%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.24'
$ws.Range("D16").Style = "Normal"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.614.16'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.06%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.049.77'
$ws.Range("D18").Style = "Normal"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -8.96%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.20'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.03%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.62'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.27%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '433.87'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.92%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.67'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -8.26%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000186'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -11.25%  '

# Row 25
$ws.Range("E25").Value = '  +5.39%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.43'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.34%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '83.94'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.58%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.56'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.08%  '

# Row 29
$ws.Range("E29").Value = '  -3.22%  '

# Row 30
$ws.Range("E30").Value = '  +0.06%  '

# Row 31
$ws.Range("E31").Value = '  +5.64%  '

# Row 32
$ws.Range("E32").Value = '  +11.98%  '

# Row 33
$ws.Range("E33").Value = '  +2.14%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.121'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -12.33%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.92'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.32%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.55'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.31%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.153'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.90%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.04'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.36%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.88'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.38%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '23.97'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.70%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.75'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.09%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.428'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.45%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '456.46'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.82%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.23'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.78%  '

# Row 46
$ws.Range("E46").Value = '  -12.07%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '160.83'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.50%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.80'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -7.77%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.655'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.05%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '43.61'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.94%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.998'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.02%  '
